# Updated cryptos list with GitHub Actions
# Applies the latest scraped coinranking.com values to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Row, $Col, $Text)
    $rng = $ws.Cells.Item($Row, $Col)
    # Force text interpretation so numeric-looking strings (e.g. "382.10",
    # "0.170", "1.85") are kept verbatim instead of being parsed as numbers,
    # then restore the default "Normal" style so no extra formatting/style
    # index is introduced on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

function Set-RowValues {
    param($Row, $B, $C, $D, $E)
    if ($null -ne $B) { Set-TextCell $Row 2 $B }
    if ($null -ne $C) { Set-TextCell $Row 3 $C }
    if ($null -ne $D) { Set-TextCell $Row 4 $D }
    if ($null -ne $E) { Set-TextCell $Row 5 $E }
}

Set-RowValues 2  $null $null "51.548.74"  "  +1.40%  "
Set-RowValues 3  $null $null "2.990.72"   "  +2.17%  "
Set-RowValues 4  $null $null $null        "  +0.05%  "
Set-RowValues 5  $null $null "382.10"     "  +1.79%  "
Set-RowValues 6  $null $null $null        "  +4.10%  "
Set-RowValues 7  $null $null $null        "  +1.82%  "
Set-RowValues 9  $null $null $null        "  +2.04%  "
Set-RowValues 10 $null $null "36.82"      "  +2.28%  "
Set-RowValues 11 $null $null $null        "  -0.60%  "
Set-RowValues 12 $null $null $null        "  +1.98%  "
Set-RowValues 13 $null $null "3.467.90"   "  +2.36%  "
Set-RowValues 14 $null $null $null        "  +3.16%  "
Set-RowValues 15 $null $null "7.79"       "  +2.69%  "
Set-RowValues 16 $null $null "2.985.91"   "  +2.20%  "
Set-RowValues 17 $null $null "11.26"      "  -0.40%  "
Set-RowValues 18 $null $null $null        "  +0.61%  "
Set-RowValues 19 $null $null "51.627.38"  "  +1.65%  "
Set-RowValues 20 $null $null $null        "  +0.77%  "
Set-RowValues 21 $null $null "12.54"      "  +1.28%  "
Set-RowValues 22 $null $null "0.0₃0966"   "  +1.30%  "
Set-RowValues 23 $null $null "70.41"      "  +2.26%  "
Set-RowValues 24 $null $null "267.70"     "  +1.03%  "
Set-RowValues 26 $null $null "8.06"       "  +0.61%  "
Set-RowValues 27 $null $null "0.170"      "  +4.71%  "
Set-RowValues 28 $null $null "7.22"       "  -2.61%  "
Set-RowValues 30 $null $null "26.12"      "  +2.45%  "
Set-RowValues 31 $null $null $null        "  +0.30%  "
Set-RowValues 32 $null $null "10.43"      "  +4.44%  "
Set-RowValues 33 $null $null "34.71"      "  +4.96%  "
Set-RowValues 34 $null $null "51.40"      "  +1.47%  "
Set-RowValues 35 $null $null $null        "  +0.58%  "
Set-RowValues 36 $null $null $null        "  +2.14%  "
Set-RowValues 37 $null $null $null        "  -0.02%  "
Set-RowValues 38 $null $null $null        "  +6.65%  "
Set-RowValues 39 $null $null "17.05"      "  +3.71%  "
Set-RowValues 40 $null $null $null        "  +5.15%  "

# Rows 41 & 42 swap places: Stellar <-> ARBITRUM
Set-RowValues 41 "ARBITRUM" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"    "1.85"  "  +2.21%  "
Set-RowValues 42 "Stellar"  "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" "0.116" "  +1.26%  "

Set-RowValues 43 $null $null "127.76"   "  +6.90%  "
Set-RowValues 44 $null $null $null      "  +14.32%  "
Set-RowValues 45 $null $null "21.47"    "  +2.24%  "
Set-RowValues 46 $null $null $null      "  +0.29%  "
Set-RowValues 47 $null $null $null      "  +1.90%  "
Set-RowValues 48 $null $null $null      "  +1.02%  "
Set-RowValues 49 $null $null "2.037.80" "  +2.51%  "
Set-RowValues 50 $null $null "3.286.79" "  +2.27%  "
Set-RowValues 51 $null $null $null      "  +2.42%  "
